$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183, pushing existing rows 183:306 down to 184:307
$ws.Rows.Item(183).Insert()

# Populate the new row 183 with the new weekly entry
$ws.Range("A183").Value = 7
$ws.Range("B183").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C183").Value = "Ñuble"
$ws.Range("D183").Value = 44767
$ws.Range("E183").Value = 16
$ws.Range("F183").Value = 100114013
$ws.Range("G183").Value = "Zanahoria"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 120
$ws.Range("K183").Value = 8000
$ws.Range("L183").Value = 9000
$ws.Range("M183").Value = 8500
$ws.Range("N183").Value = "`$/saco 20 kilos"
$ws.Range("O183").Value = "Provincia de Diguillín"
$ws.Range("P183").Value = 425
$ws.Range("Q183").Value = 20
$ws.Range("R183").Value = "Hortaliza"
